$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.369.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'1.592.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.57%  '
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("E5").Value = '  +0.94%  '
$ws.Range("E6").Value = '  +0.32%  '
$ws.Range("E7").Value = '  -0.32%  '
$ws.Range("E8").Value = '  +0.63%  '
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("E10").Value = '  -0.53%  '
$ws.Range("D11").Value = "'0.0848"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.42%  '
$ws.Range("E12").Value = '  +0.58%  '
$ws.Range("D13").Value = "'1.594.67"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.05%  '
$ws.Range("E14").Value = '  +1.23%  '
$ws.Range("E15").Value = '  +1.45%  '
$ws.Range("D16").Value = "'64.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.70%  '
$ws.Range("D17").Value = "'26.371.14"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.41%  '
$ws.Range("E18").Value = '  -1.06%  '
$ws.Range("D19").Value = "'7.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.60%  '
$ws.Range("E20").Value = '  +2.68%  '
$ws.Range("E21").Value = '  -0.34%  '
$ws.Range("E22").Value = '  +1.40%  '
$ws.Range("E23").Value = '  +2.39%  '
$ws.Range("E24").Value = '  -2.42%  '
$ws.Range("D25").Value = "'144.27"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.24%  '
$ws.Range("E26").Value = '  -0.35%  '
$ws.Range("E27").Value = '  +0.95%  '
$ws.Range("E28").Value = '  -0.45%  '
$ws.Range("D29").Value = "'15.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.03%  '
$ws.Range("D30").Value = "'0.0504"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.30%  '
$ws.Range("E31").Value = '  +1.29%  '
$ws.Range("D33").Value = "'2.98"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.17%  '
$ws.Range("D34").Value = "'1.336.99"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.12%  '
$ws.Range("E35").Value = '  -1.25%  '
$ws.Range("E36").Value = '  -0.27%  '
$ws.Range("E37").Value = '  +0.23%  '
$ws.Range("E38").Value = '  +0.21%  '
$ws.Range("E39").Value = '  +0.52%  '
$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").Value = "'1.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -19.02%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = "'5.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.56%  '
$ws.Range("E42").Value = '  -0.33%  '
$ws.Range("E43").Value = '  +0.28%  '
$ws.Range("E44").Value = '  -0.26%  '
$ws.Range("D45").Value = "'1.728.92"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.54%  '
$ws.Range("E46").Value = '  -0.64%  '
$ws.Range("D47").Value = "'88.18"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.69%  '
$ws.Range("E48").Value = '  -3.50%  '
$ws.Range("D49").Value = "'0.0987"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.98%  '
$ws.Range("E50").Value = '  -0.86%  '
$ws.Range("D51").Value = "'0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.44%  '
